# Update countries & provincias Spain
# - refresh the "last updated" timestamp
# - refresh the Covid-19 per-country stats
# - Sri Lanka overtakes Estonia in the ranking (rows 106/107 swap place)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp (cell A1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Junio de 2020 a las 19:50"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 2403956
$ws.Range("C4").Value = 15803
$ws.Range("D4").Value = 1005677
$ws.Range("E4").Value = 1275197
$ws.Range("G4").Value = 472
$ws.Range("H4").Value = 123082

# --- Row 5: Brasil ---
$ws.Range("B5").Value = 1117430
$ws.Range("C5").Value = 6082
$ws.Range("E5").Value = 471824
$ws.Range("G5").Value = 95
$ws.Range("H5").Value = 51502

# --- Row 7: India ---
$ws.Range("B7").Value = 455830
$ws.Range("C7").Value = 15380
$ws.Range("D7").Value = 258523
$ws.Range("E7").Value = 182824
$ws.Range("G7").Value = 468
$ws.Range("H7").Value = 14483

# --- Row 9: España ---
$ws.Range("B9").Value = 293832
$ws.Range("C9").Value = 248
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 28325

# --- Row 14: Alemania ---
$ws.Range("B14").Value = 192532
$ws.Range("C14").Value = 413
$ws.Range("E14").Value = 7853
$ws.Range("G14").Value = 10
$ws.Range("H14").Value = 8979

# --- Row 15: Turquia ---
$ws.Range("B15").Value = 190165
$ws.Range("C15").Value = 1268
$ws.Range("D15").Value = 162848
$ws.Range("E15").Value = 22316
$ws.Range("G15").Value = 27
$ws.Range("H15").Value = 5001

# --- Row 48: Irlanda ---
$ws.Range("B48").Value = 25391
$ws.Range("C48").Value = 8
$ws.Range("E48").Value = 973
$ws.Range("G48").Value = 3
$ws.Range("H48").Value = 1720

# --- Row 51: Israel ---
$ws.Range("B51").Value = 21467
$ws.Range("C51").Value = 385
$ws.Range("D51").Value = 15860
$ws.Range("E51").Value = 5299
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 308

# --- Row 65: Argelia ---
$ws.Range("B65").Value = 12076
$ws.Range("C65").Value = 156
$ws.Range("D65").Value = 8674
$ws.Range("E65").Value = 2541

# --- Row 68: Marruecos ---
$ws.Range("B68").Value = 10344
$ws.Range("C68").Value = 172
$ws.Range("D68").Value = 8407
$ws.Range("E68").Value = 1723

# --- Row 76: Uzbekistan ---
$ws.Range("D76").Value = 4560
$ws.Range("E76").Value = 1956

# --- Row 87: Republica de Yibuti ---
$ws.Range("B87").Value = 4617
$ws.Range("C87").Value = 18
$ws.Range("D87").Value = 3989
$ws.Range("E87").Value = 579
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 49

# --- Row 104: Maldivas ---
$ws.Range("B104").Value = 2238
$ws.Range("C104").Value = 21
$ws.Range("E104").Value = 417

# --- Rows 106/107: Sri Lanka overtakes Estonia, ranking swaps ---
# Row 106 becomes Sri Lanka (new data); row 107 becomes Estonia (prior row-106 data)
$ws.Range("A106").Value = "Sri Lanka"
$ws.Range("B106").Value = 1991
$ws.Range("C106").Value = 40
$ws.Range("D106").Value = 1548
$ws.Range("E106").Value = 432
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 11

$ws.Range("A107").Value = "Estonia"
$ws.Range("B107").Value = 1982
$ws.Range("C107").Value = 1
$ws.Range("D107").Value = 1771
$ws.Range("E107").Value = 142
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 69

# --- Row 121: Paraguay ---
$ws.Range("B121").Value = 1422
$ws.Range("C121").Value = 30
$ws.Range("D121").Value = 926
$ws.Range("E121").Value = 483

# --- Row 122: Sierra Leona ---
$ws.Range("B122").Value = 1347
$ws.Range("C122").Value = 7
$ws.Range("D122").Value = 853
$ws.Range("E122").Value = 439

# --- Row 162: Comoras ---
$ws.Range("B162").Value = 265
$ws.Range("C162").Value = 18
$ws.Range("G162").Value = 2
$ws.Range("H162").Value = 7
